$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header labels (row 1)
$ws.Range("A1").Value = "Site"
$ws.Range("B1").Value = "LCT"
$ws.Range("C1").Value = "CTI"
$ws.Range("D1").Value = "Res_Corr"
$ws.Range("E1").Value = "Aridity"
$ws.Range("F1").Value = "CTI_Class"
$ws.Range("G1").Value = "Thk"

# Add new "Thk" column values for rows 2-11
$thk = @(50, 50, 50, 0, 1, 50, 15, 39, 1, 1)
for ($i = 0; $i -lt $thk.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 7).Value = $thk[$i]
}
